$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Burndown chart 3 data fix: one day had 2 hours of work logged but 3 hours
# were actually done, so update the "Hours" value for that day from 2 to 3.
$ws.Range("E52").Value = 3

# Recalculate so the dependent SUM/burndown formulas in rows 60-61 (and the
# chart caches that reference them) pick up the new total.
$excel.Calculate()

# Update the active view/selection to match where the edit was made.
$ws.Activate()
$ws.Range("J61").Select()
$excel.ActiveWindow.ScrollRow = 46
